# Apply data corrections to "Croatia 3NL" worksheet
# - Reorders three groups of team-name entries in the shared-string table
#   (expressed here as the corrected team names written into the cells
#   that reference them), and
# - Fixes five pairs of rows whose match statistics were swapped between
#   two adjacently-id'd fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 7126859
$ws.Range("E7").Value = "NK Granicar Zupanja"
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("L7").Value = 1.5
$ws.Range("M7").Value = 4
$ws.Range("N7").Value = 5
$ws.Range("O7").Value = 1.5
$ws.Range("P7").Value = 4
$ws.Range("Q7").Value = 5
$ws.Range("R7").Value = -1
$ws.Range("S7").Value = 1.8
$ws.Range("T7").Value = 2
$ws.Range("V7").Value = 1.85
$ws.Range("W7").Value = 1.95
$ws.Range("Y7").Value = 3
$ws.Range("AA7").Value = -1
$ws.Range("AB7").Value = 1
$ws.Range("E8").Value = "NK Bistra"
$ws.Range("E9").Value = "Sava Strmec"
$ws.Range("B10").Value = 7126860
$ws.Range("E10").Value = "NK Tomislav"
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 2.5
$ws.Range("M10").Value = 3.4
$ws.Range("N10").Value = 2.4
$ws.Range("O10").Value = 2.625
$ws.Range("P10").Value = 3.4
$ws.Range("Q10").Value = 2.3
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 1.975
$ws.Range("T10").Value = 1.725
$ws.Range("V10").Value = 1.95
$ws.Range("W10").Value = 1.85
$ws.Range("Y10").Value = 2.4
$ws.Range("AA10").Value = 0
$ws.Range("AB10").Value = 0
$ws.Range("F12").Value = "Sava Strmec"
$ws.Range("F13").Value = "NK Tomislav"
$ws.Range("E19").Value = "Sava Strmec"
$ws.Range("E21").Value = "NK Bistra"
$ws.Range("F24").Value = "Sava Strmec"
$ws.Range("F27").Value = "NK Oriolik Oriovac"
$ws.Range("F28").Value = "NK Bistra"
$ws.Range("B29").Value = 7250137
$ws.Range("E29").Value = "NK Granicar Zupanja"
$ws.Range("G29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 1
$ws.Range("K29").Value = "A"
$ws.Range("L29").Value = 1.727
$ws.Range("M29").Value = 3.75
$ws.Range("N29").Value = 3.75
$ws.Range("O29").Value = 1.727
$ws.Range("P29").Value = 3.75
$ws.Range("Q29").Value = 3.75
$ws.Range("R29").Value = -0.75
$ws.Range("S29").Value = 1.975
$ws.Range("T29").Value = 1.825
$ws.Range("U29").Value = 2.5
$ws.Range("V29").Value = 1.8
$ws.Range("W29").Value = 2
$ws.Range("Y29").Value = -1
$ws.Range("Z29").Value = 2.75
$ws.Range("AA29").Value = -1
$ws.Range("AB29").Value = 0.825
$ws.Range("AC29").Value = -1
$ws.Range("B30").Value = 7250138
$ws.Range("E30").Value = "NK Tomislav"
$ws.Range("G30").Value = 2
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = "D"
$ws.Range("L30").Value = 2.1
$ws.Range("M30").Value = 3.4
$ws.Range("N30").Value = 2.9
$ws.Range("O30").Value = 2.1
$ws.Range("P30").Value = 3.4
$ws.Range("Q30").Value = 2.9
$ws.Range("R30").Value = -0.25
$ws.Range("S30").Value = 1.9
$ws.Range("T30").Value = 1.9
$ws.Range("U30").Value = 3
$ws.Range("V30").Value = 1.9
$ws.Range("W30").Value = 1.9
$ws.Range("Y30").Value = 2.4
$ws.Range("Z30").Value = -1
$ws.Range("AA30").Value = -0.5
$ws.Range("AB30").Value = 0.45
$ws.Range("AC30").Value = 0.8999999999999999
$ws.Range("F32").Value = "Slavija Pleternica"
$ws.Range("B33").Value = 7291472
$ws.Range("E33").Value = "NK Lukavec"
$ws.Range("F33").Value = "Sava Strmec"
$ws.Range("G33").Value = 3
$ws.Range("I33").Value = 1
$ws.Range("L33").Value = 2.2
$ws.Range("M33").Value = 3.6
$ws.Range("N33").Value = 2.6
$ws.Range("O33").Value = 2.2
$ws.Range("P33").Value = 3.6
$ws.Range("Q33").Value = 2.625
$ws.Range("R33").Value = -0.25
$ws.Range("S33").Value = 2
$ws.Range("T33").Value = 1.8
$ws.Range("U33").Value = 2.5
$ws.Range("V33").Value = 1.8
$ws.Range("W33").Value = 2
$ws.Range("X33").Value = 1.2
$ws.Range("AA33").Value = 1
$ws.Range("AC33").Value = 0.8
$ws.Range("B34").Value = 7291473
$ws.Range("E34").Value = "NK Maksimir"
$ws.Range("F34").Value = "NK Mladost Petrinja"
$ws.Range("G34").Value = 5
$ws.Range("I34").Value = 4
$ws.Range("L34").Value = 1.25
$ws.Range("M34").Value = 6
$ws.Range("N34").Value = 7
$ws.Range("O34").Value = 1.25
$ws.Range("P34").Value = 6
$ws.Range("Q34").Value = 7.5
$ws.Range("R34").Value = -1.75
$ws.Range("S34").Value = 1.9
$ws.Range("T34").Value = 1.9
$ws.Range("U34").Value = 3
$ws.Range("V34").Value = 1.975
$ws.Range("W34").Value = 1.825
$ws.Range("X34").Value = 0.25
$ws.Range("AA34").Value = 0.8999999999999999
$ws.Range("AC34").Value = 0.9750000000000001
$ws.Range("E38").Value = "Sava Strmec"
$ws.Range("F38").Value = "NK Bistra"
$ws.Range("E43").Value = "NK Tomislav"
$ws.Range("F45").Value = "Sloga Nova Gradiska"
$ws.Range("B50").Value = 7382546
$ws.Range("E50").Value = "NK Primorac Biograd"
$ws.Range("F50").Value = "NK Vodice"
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 1
$ws.Range("I50").Value = 1
$ws.Range("J50").Value = 1
$ws.Range("K50").Value = "D"
$ws.Range("L50").Value = 1.363
$ws.Range("M50").Value = 4.333
$ws.Range("N50").Value = 6.5
$ws.Range("O50").Value = 1.25
$ws.Range("P50").Value = 5
$ws.Range("Q50").Value = 9
$ws.Range("R50").Value = -1.75
$ws.Range("S50").Value = 1.975
$ws.Range("T50").Value = 1.825
$ws.Range("V50").Value = 1.8
$ws.Range("W50").Value = 2
$ws.Range("X50").Value = -1
$ws.Range("Y50").Value = 4
$ws.Range("AA50").Value = -1
$ws.Range("AB50").Value = 0.825
$ws.Range("AC50").Value = -1
$ws.Range("B51").Value = 7382547
$ws.Range("E51").Value = "Zmaj Makarska"
$ws.Range("F51").Value = "RNK Split"
$ws.Range("G51").Value = 5
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 5
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = "H"
$ws.Range("L51").Value = 2.3
$ws.Range("M51").Value = 3.4
$ws.Range("N51").Value = 2.625
$ws.Range("O51").Value = 2.3
$ws.Range("P51").Value = 3.4
$ws.Range("Q51").Value = 2.625
$ws.Range("R51").Value = -0.25
$ws.Range("S51").Value = 2.05
$ws.Range("T51").Value = 1.75
$ws.Range("V51").Value = 1.975
$ws.Range("W51").Value = 1.825
$ws.Range("X51").Value = 1.3
$ws.Range("Y51").Value = -1
$ws.Range("AA51").Value = 1.05
$ws.Range("AB51").Value = -1
$ws.Range("AC51").Value = 0.9750000000000001
$ws.Range("E53").Value = "NK Tomislav"
$ws.Range("E54").Value = "NK Bistra"
$ws.Range("F68").Value = "Sloga Nova Gradiska"
$ws.Range("F70").Value = "NK Bistra"
$ws.Range("F71").Value = "NK Tomislav"
$ws.Range("E78").Value = "Sava Strmec"
$ws.Range("E79").Value = "NK Bistra"
$ws.Range("F84").Value = "Sava Strmec"
$ws.Range("F88").Value = "Slavija Pleternica"
$ws.Range("E92").Value = "Sava Strmec"
$ws.Range("F94").Value = "NK Tomislav"
$ws.Range("F100").Value = "NK Oriolik Oriovac"
$ws.Range("E101").Value = "Sava Strmec"
$ws.Range("E102").Value = "NK Bistra"
$ws.Range("B109").Value = 8061519
$ws.Range("E109").Value = "NK Bistra"
$ws.Range("F109").Value = "NK Lukavec"
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 1
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 2.5
$ws.Range("M109").Value = 3.4
$ws.Range("N109").Value = 2.375
$ws.Range("O109").Value = 2.5
$ws.Range("P109").Value = 3.4
$ws.Range("Q109").Value = 2.375
$ws.Range("R109").Value = 0
$ws.Range("S109").Value = 1.975
$ws.Range("T109").Value = 1.825
$ws.Range("U109").Value = 2.75
$ws.Range("V109").Value = 1.825
$ws.Range("W109").Value = 1.975
$ws.Range("Z109").Value = 1.375
$ws.Range("AB109").Value = 0.825
$ws.Range("AC109").Value = -1
$ws.Range("B110").Value = 8061520
$ws.Range("F110").Value = "HNK Segesta"
$ws.Range("G110").Value = 2
$ws.Range("H110").Value = 6
$ws.Range("I110").Value = 1
$ws.Range("J110").Value = 4
$ws.Range("L110").Value = 5
$ws.Range("M110").Value = 4
$ws.Range("N110").Value = 1.5
$ws.Range("O110").Value = 6.5
$ws.Range("P110").Value = 4.2
$ws.Range("Q110").Value = 1.4
$ws.Range("R110").Value = 1.25
$ws.Range("S110").Value = 1.825
$ws.Range("T110").Value = 1.975
$ws.Range("U110").Value = 2.5
$ws.Range("V110").Value = 1.8
$ws.Range("W110").Value = 2
$ws.Range("Z110").Value = 0.3999999999999999
$ws.Range("AB110").Value = 0.9750000000000001
$ws.Range("AC110").Value = 0.8
$ws.Range("F119").Value = "Sloga Nova Gradiska"
$ws.Range("E120").Value = "NK Bistra"
$ws.Range("F120").Value = "Sava Strmec"
$ws.Range("F122").Value = "Slavija Pleternica"
$ws.Range("E125").Value = "Sava Strmec"
$ws.Range("E128").Value = "NK Bistra"
$ws.Range("F129").Value = "Sava Strmec"
$ws.Range("E134").Value = "NK Bistra"
$ws.Range("E141").Value = "Sava Strmec"
$ws.Range("E142").Value = "NK Bistra"
